# Insert a new weekly price record as row 110 (Feria Lagunitas de Puerto Montt -
# Acelga), pushing the existing rows 110-131 down to 111-132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 110; Excel shifts rows 110..131 -> 111..132
# and carries the column-D date formatting down from the row above, same as a
# manual "Insert Sheet Rows" in the UI.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record.
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 44551
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = 100112009
$ws.Cells.Item(110, 7).Value = "Acelga"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 50
$ws.Cells.Item(110, 11).Value = 10000
$ws.Cells.Item(110, 12).Value = 10000
$ws.Cells.Item(110, 13).Value = 10000
$ws.Cells.Item(110, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(110, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(110, 16).Value = 833
$ws.Cells.Item(110, 17).Value = 12
$ws.Cells.Item(110, 18).Value = "Hortaliza"
